# Updated cryptos list on Wed May 22 11:32:02 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, [string]$addr, [string]$val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "70.071.20"
Set-TextValue $ws "E2" "  -1.60%  "

Set-TextValue $ws "D3" "3.733.51"
Set-TextValue $ws "E3" "  -0.89%  "

Set-TextValue $ws "E4" "  +0.07%  "

Set-TextValue $ws "D5" "621.03"
Set-TextValue $ws "E5" "  +0.20%  "

Set-TextValue $ws "D6" "180.50"
Set-TextValue $ws "E6" "  -1.04%  "

Set-TextValue $ws "D7" "3.731.99"
Set-TextValue $ws "E7" "  -0.74%  "

Set-TextValue $ws "E8" "  +0.02%  "

Set-TextValue $ws "D9" "0.534"
Set-TextValue $ws "E9" "  -1.67%  "

Set-TextValue $ws "D10" "0.167"
Set-TextValue $ws "E10" "  +1.35%  "

Set-TextValue $ws "D11" "6.30"
Set-TextValue $ws "E11" "  -4.75%  "

Set-TextValue $ws "D12" "0.485"
Set-TextValue $ws "E12" "  -3.61%  "

Set-TextValue $ws "D13" "40.65"
Set-TextValue $ws "E13" "  -0.07%  "

Set-TextValue $ws "D14" "0.0000258"
Set-TextValue $ws "E14" "  +0.53%  "

Set-TextValue $ws "D15" "4.352.32"
Set-TextValue $ws "E15" "  -0.68%  "

Set-TextValue $ws "D16" "3.733.23"
Set-TextValue $ws "E16" "  -0.69%  "

Set-TextValue $ws "D17" "70.110.08"
Set-TextValue $ws "E17" "  -1.62%  "

Set-TextValue $ws "E18" "  -1.69%  "

Set-TextValue $ws "E19" "  +0.34%  "

Set-TextValue $ws "D20" "16.80"
Set-TextValue $ws "E20" "  -0.87%  "

Set-TextValue $ws "D21" "505.92"
Set-TextValue $ws "E21" "  -2.68%  "

Set-TextValue $ws "D22" "9.30"
Set-TextValue $ws "E22" "  -0.93%  "

Set-TextValue $ws "D23" "0.723"
Set-TextValue $ws "E23" "  -3.54%  "

Set-TextValue $ws "E24" "  +0.50%  "

Set-TextValue $ws "D25" "86.84"
Set-TextValue $ws "E25" "  -2.32%  "

Set-TextValue $ws "D26" "11.49"
Set-TextValue $ws "E26" "  +2.13%  "

Set-TextValue $ws "D27" "13.11"
Set-TextValue $ws "E27" "  -3.75%  "

Set-TextValue $ws "E28" "  +21.28%  "

Set-TextValue $ws "E29" "  -0.23%  "

Set-TextValue $ws "E30" "  -2.44%  "

Set-TextValue $ws "D31" "2.93"
Set-TextValue $ws "E31" "  +0.79%  "

Set-TextValue $ws "D32" "7.93"
Set-TextValue $ws "E32" "  -3.03%  "

Set-TextValue $ws "D33" "31.19"
Set-TextValue $ws "E33" "  -3.60%  "

Set-TextValue $ws "E34" "  -1.15%  "

Set-TextValue $ws "D35" "1.00"
Set-TextValue $ws "E35" "  +0.25%  "

Set-TextValue $ws "E36" "  +0.42%  "

Set-TextValue $ws "E37" "  -0.04%  "

Set-TextValue $ws "E38" "  +2.00%  "

Set-TextValue $ws "D39" "0.340"
Set-TextValue $ws "E39" "  -1.78%  "

Set-TextValue $ws "D40" "2.08"
Set-TextValue $ws "E40" "  -6.78%  "

Set-TextValue $ws "D41" "50.34"
Set-TextValue $ws "E41" "  -2.87%  "

Set-TextValue $ws "D42" "45.58"
Set-TextValue $ws "E42" "  +1.58%  "

Set-TextValue $ws "D43" "433.95"
Set-TextValue $ws "E43" "  -1.16%  "

Set-TextValue $ws "D44" "2.89"
Set-TextValue $ws "E44" "  +1.64%  "

Set-TextValue $ws "D45" "8.70"
Set-TextValue $ws "E45" "  -2.00%  "

Set-TextValue $ws "D46" "3.007.91"
Set-TextValue $ws "E46" "  -5.38%  "

Set-TextValue $ws "D47" "0.0363"
Set-TextValue $ws "E47" "  -0.73%  "

Set-TextValue $ws "D48" "27.48"
Set-TextValue $ws "E48" "  -2.59%  "

Set-TextValue $ws "D50" "137.53"
Set-TextValue $ws "E50" "  -2.26%  "

Set-TextValue $ws "E51" "  +0.04%  "
